# This edit shuffles the per-observation records stored in rows 15-38 of the
# sheet (one row per species finding) among themselves: every row in that
# block (except row 17, which keeps its own data) ends up holding the data
# that originally belonged to a *different* row in the same block. No rows
# are inserted/removed and no other part of the sheet changes.
#
# Strategy: capture every cell (columns A:AY) of every affected row from the
# sheet exactly as it is *before* any writes happen, then write the captured
# values into their new (target) rows. Capturing everything up front avoids
# any ordering problems from the permutation overwriting source data that is
# still needed for a later target row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows holding the shuffled records.
$rows = @(15,16,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38)
$maxCol = 51  # column AY is the right edge of the used range

# Step 1: snapshot every cell in the affected rows before making any changes.
$data = @{}
foreach ($r in $rows) {
    for ($c = 1; $c -le $maxCol; $c++) {
        $key = "$r" + "_" + "$c"
        $data[$key] = $ws.Cells.Item($r, $c).Value()
    }
}

# Step 2: target row -> source row. The target row receives the snapshot
# taken (in step 1) from the source row.
$map = @{}
$map[15] = 20
$map[16] = 32
$map[18] = 19
$map[19] = 27
$map[20] = 33
$map[21] = 23
$map[22] = 35
$map[23] = 22
$map[24] = 37
$map[25] = 29
$map[26] = 18
$map[27] = 30
$map[28] = 21
$map[29] = 26
$map[30] = 28
$map[31] = 38
$map[32] = 34
$map[33] = 16
$map[34] = 15
$map[35] = 31
$map[36] = 24
$map[37] = 36
$map[38] = 25

# Step 3: write the snapshots into their new rows.
foreach ($t in $map.Keys) {
    $s = $map[$t]
    for ($c = 1; $c -le $maxCol; $c++) {
        $key = "$s" + "_" + "$c"
        $v = $data[$key]
        # The date string "2023-08-24" (columns Y/AA) would otherwise be
        # auto-detected as a real date by the normal Value setter and turned
        # into a date serial number; a leading apostrophe forces it to stay
        # literal text, matching the source data (which stores it as text).
        if (($v -is [string]) -and ($v -eq "2023-08-24")) {
            $ws.Cells.Item($t, $c).Value = "'" + $v
        } else {
            $ws.Cells.Item($t, $c).Value = $v
        }
    }
}
